# Database_Thresholds.xlsx update: bump script-run date, update quantile-source
# commit hash, add a "Script Run" header row, and fix two data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (rows 1-3) ---------------------------------------------
# Row 1 = title (untouched). Row 2 is a new line carrying the refreshed
# "Script Run" stamp (previously this text lived in row 3). Row 3 keeps its
# bold styling but becomes blank.
$ws.Cells.Item(2, 1).Value2 = "Script Run: 2024-11-25"
$ws.Cells.Item(3, 1).Value2 = ""

# --- Per-row script metadata (rows 8-96: ScriptLatestRunDate in col X) ---
for ($r = 8; $r -le 96; $r++) {
    $ws.Cells.Item($r, 24).Value2 = 45621
}

# --- Rows that also carry ActionNeededDate / QuantileSource / QuantileDate
$actionRows = @(67,68,69,70,73,74,75,76,77,78,79,82,83,84,85,87,88,89,91,92,94)
foreach ($r in $actionRows) {
    $ws.Cells.Item($r, 20).Value2 = 45621
    $ws.Cells.Item($r, 21).Value2 = "Database_Thresholds.xlsx, Git Commit ID: 00a27356deb3dd3ebc41bb112fe7b3bfe3e44544"
    $ws.Cells.Item($r, 25).Value2 = 45621
}

# --- Two corrected data values --------------------------------------------
$ws.Cells.Item(68, 15).Value2 = 2.837179
$ws.Cells.Item(87, 10).Value2 = 25
